# Insert a new data row at row 112 (pushing existing rows 112..158 down to 113..159)
# and populate it with the new price-listing record for Piña (pineapple), matching
# the commit "Fruta / hortaliza, semanal" weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("112:112").Insert()

$ws.Range("A112").Value = 5
$ws.Range("B112").Value = 'Macroferia Regional de Talca'
$ws.Range("C112").Value = 'Maule'
$ws.Range("D112").Value = 44489
$ws.Range("E112").Value = 7
$ws.Range("F112").Value = 'Fruta'
$ws.Range("G112").Value = 100108
$ws.Range("H112").Value = 'Tropicales y subtropicales'
$ws.Range("I112").Value = 100108005
$ws.Range("J112").Value = 'Piña'
$ws.Range("K112").Value = 'Caramelo'
$ws.Range("L112").Value = 'Segunda'
$ws.Range("M112").Value = 200
$ws.Range("N112").Value = 19000
$ws.Range("O112").Value = 19000
$ws.Range("P112").Value = 19000
$ws.Range("Q112").Value = '$/caja 14 unidades'
$ws.Range("R112").Value = 'Ecuador'
$ws.Range("S112").Value = 1357
$ws.Range("T112").Value = 14
